# "Typology using UPGMA done"
# - Renames the placeholder "Sheet2" to "groups_UPGMA"
# - Fills in the "broader typology_UPGMA" cluster-description table
# - Updates sheet selections / active tab to match the saved state

$wb = $excel.ActiveWorkbook

# 1) Rename the still-blank UPGMA groups sheet
$wsGroups = $wb.Worksheets.Item("Sheet2")
$wsGroups.Name = "groups_UPGMA"

# Leave the cursor where the author left it on that sheet
$wsGroups.Range("R6").Select()

# 2) Populate the "broader typology_UPGMA" summary table
$wsTypo = $wb.Worksheets.Item("broader typology_UPGMA")

# Widen the description column so the paragraphs are readable
$wsTypo.Columns.Item(2).ColumnWidth = 127.25

# Cluster 1 (row 2)
$wsTypo.Range("B2").Value = "Very large provinces at very high elevations. Very low education levels, very high proportion of primary sector workers and very low proportion of scondary sector workers. Economic security provided by rural livelihoods - few people have no farmland and livestock ownership is common. Very low access to services, high crime per capita, low land conflict and very low migration levels."
$wsTypo.Range("C2").Value = "VERY HIGH"

# Cluster 2 (row 3)
$wsTypo.Range("B3").Value = "Very small provinces at very high elevations. Low levels of education, low proportion of people in the primary sector but higher proportion of people in the secondary sector. Very few people with no farmland, but very little livestock ownership. High access to services and high crime per capita. Low land conflict and low migration."
$wsTypo.Range("C3").Value = "LOW"

# Cluster 3 (row 4)
$wsTypo.Range("B4").Value = "Small provinces at very low elevations. Very high levels of education, high proportion of people in the primary sector, but very high proportion of people in the secondary sector. High proportion of people with no farmland, but high levels of livestock ownership. High access to services and low crime per capita. But very high migration levels and very high rates of land conflict."
$wsTypo.Range("C4").Value = "VERY LOW"

# Cluster 4 (row 5)
$wsTypo.Range("B5").Value = "Large provinces at low elevations. Very high levels of education, and relatively low proportion of workers in the primary and secondary sectors (suggesting higher proportions in the other sectors e.g. tertiary). High proportion of people with no farmland, and low levels of livestock ownership (suggesting very urban).  Low access to services, but this may be explained by the mean size of the provinces in this cluster (there is high access to garbage collection). Low crime per capita, but very high migration and very high rates of land conflict"
$wsTypo.Range("C5").Value = "HIGH"

# Cluster 5 (row 6)
$wsTypo.Range("B6").Value = "Very large provinces at high elevations. High levels of education, and a high proportion of workers in both primary and secondary sectors. Very high proportion of people with no farmland, but also very high proportion of people with livestock. Low access to services (although very high access to garbage collection) - this may be an artefact of the very large mean area of the provinces in this cluster. Very high crime rates, very high migration, and very high rates of land conflict."
$wsTypo.Range("C6").Value = "VERY HIGH"

# Wrap the long description cells like the sibling "Broader typology_kmeans" sheet
$wsTypo.Range("B2:B6").WrapText = $true

# Row heights for the wrapped paragraphs
$wsTypo.Rows.Item(2).RowHeight = 43.2
$wsTypo.Rows.Item(3).RowHeight = 43.2
$wsTypo.Rows.Item(4).RowHeight = 43.2
$wsTypo.Rows.Item(5).RowHeight = 57.6
$wsTypo.Rows.Item(6).RowHeight = 57.6

# This is the sheet the author left open/active when they saved
$wsTypo.Range("B7").Select()
